# Applies the "minor changes to a few PowerPoint slides" edit:
#   - Slide 13: reword two bullet points about the course assembler.
#   - Slide 17: "x86" -> "x86-64" label (inside a grouped shape).
#   - Slide 19: reword a bullet point; reposition a caption textbox.
#   - Slide 20: reword "object code" -> "assembly code" (and keep the
#     auto-fit text box's height unchanged, matching the source diff).
#   - Slide 8: reposition a caption textbox (no text change).

$p = $ppt.ActivePresentation

# --- Slide 13: "Final Code Generator" -------------------------------------
$s13 = $p.Slides.Item(13)
$body13 = $s13.Shapes.Item(4)
$body13.TextFrame.TextRange.Paragraphs(4).Runs(1).Text = "assembler provided in the course repository performs final code generation"
$body13.TextFrame.TextRange.Paragraphs(5).Runs(1).Text = "assembler also implements minor optimizations"

# --- Slide 17: "Compiler Construction Tools" diagram -----------------------
$s17 = $p.Slides.Item(17)
$group17 = $s17.Shapes.Item(5)
$x86Box = $group17.GroupItems.Item(4)
$x86Box.TextFrame.TextRange.Paragraphs(1).Runs(1).Text = "x86-64"

# --- Slide 19: "Single-pass Versus Multi-pass Compilers" -------------------
$s19 = $p.Slides.Item(19)
$body19 = $s19.Shapes.Item(4)
$body19.TextFrame.TextRange.Paragraphs(10).Runs(1).Text = "requires design of intermediate languages/representations"

$caption19 = $s19.Shapes.Item(5)
$caption19.Top = 450.0

# --- Slide 20: "Passes in the Compiler Project" -----------------------------
$s20 = $p.Slides.Item(20)
$note20 = $s20.Shapes.Item(5)
$originalHeight20 = $note20.Height
$note20.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "in-memory data structures called abstract syntax trees.  The only I/O to disk occurs when reading the source file and generating assembly code."
# Re-assert the box height: editing the text re-flows the auto-fit shape,
# but the source deck keeps the original size, so restore it explicitly.
$note20.Height = $originalHeight20

# --- Slide 8: "Constraint Analyzer" -----------------------------------------
$s8 = $p.Slides.Item(8)
$caption8 = $s8.Shapes.Item(5)
$caption8.Top = 384.5671653543307
